# Update computation of adverse event costs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (column H is new) ---
$ws.Range("H1").Value = "notes"

# --- Column widths for new G/H note columns ---
# (COM ColumnWidth snaps to a fixed 1/6-character grid; these inputs land on
# the closest achievable grid point to the target widths 41.1640625 / 43.6640625)
$ws.Range("G1").EntireColumn.ColumnWidth = 40.25
$ws.Range("H1").EntireColumn.ColumnWidth = 42.75

# --- Row 2 (Diarrhea): reference changes from "DRG 808" to "DRG 391" ---
$ws.Range("G2").Value = "DRG 391"
$ws.Range("H2").Value = ""

# --- Row 3 (Dry skin): cost value + reference + note ---
$ws.Range("C3").Value = 940
$ws.Range("G3").Value = 'Wong, William, et al. "Assessment of costs associated with adverse events in patients with cancer." PloS one 13.4 (2018): e0196007.'
$ws.Range("H3").Value = "Used the cost estimate for dermatitis"

# --- Row 4 (Elevated alanine transaminase): cost value + reference + note ---
$ws.Range("C4").Value = 3900
$ws.Range("C4").NumberFormat = "#,##0"
$ws.Range("G4").Value = 'Latremouille-Viau, Dominick, et al. "The economic burden of common adverse events associated with metastatic colorectal cancer treatment in the United States." Journal of medical economics 20.1 (2017): 54-62.'
$ws.Range("H4").Value = "Used the adjusted monthly cost difference b/w individuals with hepatobiliary/pancreatic AEs and those without and multiplied by 12 for yearly estimate"

# --- Row 5 (Elevated aspartate transaminase): cost value + reference + note ---
$ws.Range("C5").Value = 3900
$ws.Range("C5").NumberFormat = "#,##0"
$ws.Range("G5").Value = 'Latremouille-Viau, Dominick, et al. "The economic burden of common adverse events associated with metastatic colorectal cancer treatment in the United States." Journal of medical economics 20.1 (2017): 54-62.'
$ws.Range("H5").Value = "Used the adjusted monthly cost difference b/w individuals with hepatobiliary/pancreatic AEs and those without and multiplied by 12 for yearly estimate"

# --- Row 6 (Eye problems): cost value + reference ---
$ws.Range("C6").Value = 2737
$ws.Range("C6").NumberFormat = "#,##0"
$ws.Range("G6").Value = 'Wong, William, et al. "Assessment of costs associated with adverse events in patients with cancer." PloS one 13.4 (2018): e0196007.'
$ws.Range("H6").Value = ""

# --- Row 7 (Paronychia): add blank styled notes cell (reference unchanged) ---
$ws.Range("H7").Value = ""

# --- Row 8 (Pneumonitis): add blank styled notes cell (reference unchanged) ---
$ws.Range("H8").Value = ""

# --- Row 9 (Pruritus): cost value + reference ---
$ws.Range("C9").Value = 1184
$ws.Range("C9").NumberFormat = "#,##0"
$ws.Range("G9").Value = 'Wong, William, et al. "Assessment of costs associated with adverse events in patients with cancer." PloS one 13.4 (2018): e0196007.'
$ws.Range("H9").Value = ""

# --- Row 10 (Rash): cost value changes, reference changes to Wong citation ---
$ws.Range("C10").Value = 940
$ws.Range("G10").Value = 'Wong, William, et al. "Assessment of costs associated with adverse events in patients with cancer." PloS one 13.4 (2018): e0196007.'
$ws.Range("H10").Value = ""

# --- Row 11 (Stomatitis): add blank styled notes cell (reference unchanged) ---
$ws.Range("H11").Value = ""

# --- Wrap text for all G/H note cells (rows 2-11) ---
$ws.Range("G2:H11").WrapText = $true

# --- Row heights for rows with long wrapped reference text ---
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 75
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 45

# --- Stray formatted empty cell at C14 (keeps the #,##0 number style) ---
$ws.Range("C14").NumberFormat = "#,##0"

# --- Selection state ---
$ws.Range("D5").Select()

Write-Host "edit complete"
